# Adds a new "optional" column (I) to the "fields" sheet.
# Header "optional" in I1, boolean values in I2:I172.
# Rows 17-41 (Tag1..Tag25 fields) are optional = TRUE; all other rows = FALSE.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("fields")

# Header cell, matching bold style of neighboring header cells (G1/H1).
$ws.Range("I1").Value = "optional"
$ws.Range("I1").Font.Bold = $true

$optionalRows = @(17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41)

for ($r = 2; $r -le 172; $r++) {
    $cell = $ws.Cells.Item($r, 9)
    if ($optionalRows -contains $r) {
        $cell.Value = $true
    } else {
        $cell.Value = $false
    }
}

# Refresh the AutoFilter to cover the new column.
$ws.AutoFilterMode = $false
[void]$ws.Range("A1:I172").AutoFilter()

# Update the hidden _FilterDatabase defined name to the new range.
$names = $wb.Names
for ($i = 1; $i -le $names.Count; $i++) {
    $n = $names.Item($i)
    if ($n.Name -eq "fields!_FilterDatabase") {
        $n.RefersTo = "=fields!`$A`$1:`$I`$172"
    }
}

# Match the recorded selection after the edit.
[void]$ws.Activate()
[void]$ws.Range("I17").Select()
